$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.814.15'
$ws.Range("E2").Value = '  -1.29%  '
$ws.Range("D3").Value = '1.856.85'
$ws.Range("E3").Value = '  -0.58%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.22'
$ws.Range("E5").Value = '  -0.92%  '
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5035'
$ws.Range("E7").Value = '  -2.78%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3653'
$ws.Range("E8").Value = '  -2.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07165'
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8914'
$ws.Range("E10").Value = '  +0.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.63'
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07510'
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("D13").Value = '1.855.01'
$ws.Range("E13").Value = '  -0.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.24'
$ws.Range("E14").Value = '  +3.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.228'
$ws.Range("E15").Value = '  -2.07%  '
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008493'
$ws.Range("E17").Value = '  -0.83%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.04'
$ws.Range("E18").Value = '  -0.81%  '
$ws.Range("D20").Value = '26.865.18'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.029'
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").Value = '2.087.89'
$ws.Range("E22").Value = '  -0.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.33'
$ws.Range("E23").Value = '  -2.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.452'
$ws.Range("E24").Value = '  -0.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.37'
$ws.Range("E25").Value = '  -2.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.794'
$ws.Range("E26").Value = '  -2.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.79'
$ws.Range("E27").Value = '  -1.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.059'
$ws.Range("E28").Value = '  -3.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.80'
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.631'
$ws.Range("E30").Value = '  -2.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.659'
$ws.Range("E31").Value = '  -0.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09210'
$ws.Range("E32").Value = '  +2.24%  '
$ws.Range("E33").Value = '  -1.40%  '
$ws.Range("E34").Value = '  -3.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7406'
$ws.Range("E35").Value = '  -1.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.145'
$ws.Range("E36").Value = '  -2.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.243'
$ws.Range("E37").Value = '  +7.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.509'
$ws.Range("E38").Value = '  -1.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01987'
$ws.Range("E39").Value = '  -2.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.082'
$ws.Range("E40").Value = '  +0.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5320'
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '119.44'
$ws.Range("E42").Value = '  +3.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.485'
$ws.Range("E43").Value = '  -2.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.372'
$ws.Range("E44").Value = '  -1.37%  '
$ws.Range("E45").Value = '  -1.66%  '
$ws.Range("E46").Value = '  -0.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9999'
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.990'
$ws.Range("E48").Value = '  -1.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.558'
$ws.Range("E49").Value = '  -0.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.91'
$ws.Range("E50").Value = '  +1.27%  '
$ws.Range("E51").Value = '  -3.33%  '
